# progression.xlsx -- "updated version of database"
#
# Logical changes (derived from the target OOXML diff):
#  - B3  : "How much f changes as x changes"
#          -> "How much $f$ changes as $x$ changes"
#  - B38 : "The components of d \vec r is an arbitrary small change ..."
#          -> "The components of $d \vec r$ is an arbitrary small change ..."
#  - B39 : "The magnitude of dr is the length of a small step along a path"
#          -> "The magnitude of $d\vec r$ is the length of a small step along a path"
#  - G39 : (empty) -> "PH 422"                                   (new cell)
#  - B40 : "The direction of dr is ? "
#          -> "The direction of $d\vec r$ is ? "
#  - G40 : (empty) -> "PH 422"                                   (new cell)
#  - B49 : "Differential form of r in spherical and cylindrical coordinates"
#          -> "Differential form of $\vec r$ in spherical and cylindrical coordinates"
#
# (Every other line touched by the raw XML diff is a mechanical side effect
# of the shared-string table being rebuilt/reordered by these edits -- the
# actual text shown to the user, e.g. all the "PH 422" cells in column G for
# rows 49-68, does not change.)
#
# Plus a couple of cosmetic/view changes also present in the diff:
#  - Row 3 grows from one line tall to two lines tall (15.75 -> 31.5) because
#    the new text wraps.
#  - The frozen-pane view scrolls back to the top-left and the selected cell
#    moves from K47 to B4.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("B3").Value  = "How much `$f`$ changes as `$x`$ changes"
$ws.Range("B38").Value = "The components of `$d \vec r`$ is an arbitrary small change between two arbitrary position vectors. "
$ws.Range("B39").Value = "The magnitude of `$d\vec r`$ is the length of a small step along a path"
$ws.Range("G39").Value = "PH 422"
$ws.Range("B40").Value = "The direction of `$d\vec r`$ is ? "
$ws.Range("G40").Value = "PH 422"
$ws.Range("B49").Value = "Differential form of `$\vec r`$ in spherical and cylindrical coordinates"

# Row 3 now wraps onto a second line.
$ws.Rows.Item(3).RowHeight = 31.5

# Restore the view to the top of the frozen pane with B4 selected.
$ws.Range("B4").Select() | Out-Null
